# Update the workbook's recorded absolute path (cosmetic file-system metadata)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new "R_spec" row before row 28 (shifts cp/kappa/.../tau block down by 1) ---
$ws.Rows(28).Insert()
$ws.Range("A28").Value = "R_spec"
$ws.Range("B28").Value = 287.058
$ws.Range("C28").Value = "[J/K kg]"

# --- Insert new "viscosity_air" row before (current) row 34, i.e. right before rho_w ---
$ws.Rows(34).Insert()
$ws.Range("A34").Value = "viscosity_air"
$ws.Range("B34").Value = 0.0000173
$ws.Range("B34").NumberFormat = "0.00E+00"
$ws.Range("C34").Value = "[kg/(m sec)]"

# --- Append two new rows at the bottom of the table ---
$ws.Range("A56").Value = "Darcy_friction_factor"
$ws.Range("B56").Value = 0.1
$ws.Range("C56").Value = "rough-pipe regime"

$ws.Range("A57").Value = "tortuosity_air"
$ws.Range("B57").Value = 2.5
$ws.Range("C57").Value = "used in Carman Kozeny model"

# --- Update the visible selection to match the saved view state ---
$ws.Range("A23").Select()
$ws.Range("C30").Select()
